$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1801.6666
$ws.Range("I8").Value = 162.2
$ws.Range("K8").Value = 486.6
$ws.Range("M8").Value = -347.6
$ws.Range("H11").Value = 34.666668
$ws.Range("I11").Value = 34.666668
$ws.Range("K11").Value = 34.666668
$ws.Range("M11").Value = 105.333332
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H28").Value = 2937.875
$ws.Range("I28").Value = 1900.8
$ws.Range("K28").Value = 1900.8
$ws.Range("M28").Value = -1415.8
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H33").Value = 72.59999999999999
$ws.Range("I33").Value = 72.59999999999999
$ws.Range("K33").Value = 72.59999999999999
$ws.Range("M33").Value = 156.4
$ws.Range("H41").Value = 453
$ws.Range("I41").Value = 335.8
$ws.Range("J41").Value = 746
$ws.Range("K41").Value = 335.8
$ws.Range("L41").Value = 746
$ws.Range("M41").Value = 104.2
$ws.Range("N41").Value = -1626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 745.4286
$ws.Range("I32").Value = 572
$ws.Range("K32").Value = 572
$ws.Range("M32").Value = -285
$ws.Range("H45").Value = 14570
$ws.Range("I45").Value = 8997.5
$ws.Range("K45").Value = 8997.5
$ws.Range("M45").Value = -8620.5
$ws.Range("H61").Value = 8588.556
$ws.Range("I61").Value = 5382.8335
$ws.Range("K61").Value = 5382.8335
$ws.Range("M61").Value = -5170.8335
$ws.Range("H132").Value = 13728.167
$ws.Range("I132").Value = 11676.857
$ws.Range("K132").Value = 35030.571
$ws.Range("M132").Value = -32500.571
$ws.Range("H136").Value = 8588.556
$ws.Range("I136").Value = 5382.8335
$ws.Range("K136").Value = 16148.5005
$ws.Range("M136").Value = -13598.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H64").Value = 3313.6667
$ws.Range("I64").Value = 1620.75
$ws.Range("J64").Value = 6699.5
$ws.Range("K64").Value = 1620.75
$ws.Range("L64").Value = 6699.5
$ws.Range("M64").Value = -1395.75
$ws.Range("N64").Value = -7149.5
$ws.Range("H67").Value = 3313.6667
$ws.Range("I67").Value = 1620.75
$ws.Range("J67").Value = 6699.5
$ws.Range("K67").Value = 1620.75
$ws.Range("L67").Value = 6699.5
$ws.Range("M67").Value = -840.75
$ws.Range("N67").Value = -8259.5
$ws.Range("H100").Value = 10174.25
$ws.Range("J100").Value = 10174.25
$ws.Range("L100").Value = 10174.25
$ws.Range("N100").Value = -12338.25
$ws.Range("H134").Value = 8424.333000000001
$ws.Range("I134").Value = 3469.8333
$ws.Range("K134").Value = 10409.4999
$ws.Range("M134").Value = -7874.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 109.6
$ws.Range("I19").Value = 109.6
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 109.6
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 199
$ws.Range("I22").Value = 199
$ws.Range("K22").Value = 199
$ws.Range("M22").Value = 151
$ws.Range("H24").Value = 109.6
$ws.Range("I24").Value = 109.6
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 109.6
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("H132").Value = 9835
$ws.Range("I132").Value = 8502.5
$ws.Range("K132").Value = 25507.5
$ws.Range("M132").Value = -22977.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 37.11111
$ws.Range("I7").Value = 37.5
$ws.Range("K7").Value = 112.5
$ws.Range("M7").Value = -0.5
$ws.Range("H13").Value = 407.2857
$ws.Range("J13").Value = 1270
$ws.Range("L13").Value = 3810
$ws.Range("N13").Value = -4146
$ws.Range("H19").Value = 8749.75
$ws.Range("J19").Value = 10000
$ws.Range("L19").Value = 30000
$ws.Range("N19").Value = -30348
$ws.Range("H121").Value = 2248.5
$ws.Range("I121").Value = 499
$ws.Range("K121").Value = 1497
$ws.Range("M121").Value = -187

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1002
$ws.Range("I7").Value = 1002
$ws.Range("K7").Value = 1002
$ws.Range("M7").Value = -890
$ws.Range("H8").Value = 1002
$ws.Range("I8").Value = 1002
$ws.Range("K8").Value = 1002
$ws.Range("M8").Value = -863
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").ClearContents()
$ws.Range("H48").Value = 30027
$ws.Range("I48").Value = 30027
$ws.Range("K48").Value = 30027
$ws.Range("M48").Value = -29542
$ws.Range("H52").Value = 40000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 40000
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -40518
$ws.Range("H80").Value = 2625
$ws.Range("J80").Value = 2833.3333
$ws.Range("L80").Value = 2833.3333
$ws.Range("N80").Value = -4829.3333
$ws.Range("H83").Value = 2625
$ws.Range("J83").Value = 2833.3333
$ws.Range("L83").Value = 14166.6665
$ws.Range("N83").Value = -24150.6665
$ws.Range("H102").Value = 4312.25
$ws.Range("I102").Value = 4312.25
$ws.Range("K102").Value = 4312.25
$ws.Range("M102").Value = -2690.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H136").Value = 10081.75
$ws.Range("I136").Value = 9098.1
$ws.Range("K136").Value = 27294.3
$ws.Range("M136").Value = -24744.3
